$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the existing header style (from H1, style index 1) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-14 for columns I (I0) and J (IF)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 5

$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 9

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 4

$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 3

$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 5

$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 5

$ws.Range("I10").Value = 5
$ws.Range("J10").Value = 8

$ws.Range("I11").Value = 8
$ws.Range("J11").Value = 9

$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 7

$ws.Range("I13").Value = 3
$ws.Range("J13").Value = 4

$ws.Range("I14").Value = 3
$ws.Range("J14").Value = 4
